{"js": "// Resume formatting fixes:\n// 1. \"2 years of experience in C or C++\" -> \"2 years of experience in C/C++\"\n// 2. \"Full Phase \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\" -> \"Full Phases \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\"\n// 3. First \"Use project tracking tool \u2013 Redmine, Jira\" gets a trailing \".\"\n// 4. First \"Experience in QA/QC workflow\" gets a trailing \".\"\n\nconst EN_DASH = \"\\u2013\";\n\n// 1) \"C or C++\" -> \"C/C++\" (the only occurrence of this sentence).\nconst cOrCpp = context.document.body.search(\"2 years of experience in C or C++\", { matchCase: true });\ncOrCpp.load(\"text\");\nawait context.sync();\nif (cOrCpp.items.length > 0) {\n  cOrCpp.items[0].insertText(\"2 years of experience in C/C++\", Word.InsertLocation.replace);\n}\n\n// 2) \"Full Phase\" -> \"Full Phases\" (the only instance of this project-name line).\nconst fullPhase = context.document.body.search(\n  \"Full Phase \" + EN_DASH + \" MCAL drivers \" + EN_DASH + \" AUTOSAR (RVC)\",\n  { matchCase: true }\n);\nfullPhase.load(\"text\");\nawait context.sync();\nif (fullPhase.items.length > 0) {\n  fullPhase.items[0].insertText(\n    \"Full Phases \" + EN_DASH + \" MCAL drivers \" + EN_DASH + \" AUTOSAR (RVC)\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 3) Append \".\" right after the FIRST \"Use project tracking tool \u2013 Redmine, Jira\"\n//    (this sentence repeats in several project bullet lists; only the first gets the\n//    trailing period). Insert at the end of its paragraph (not \"after\" the found\n//    range) so the existing \"-\" / tab runs that precede it in the bullet line are\n//    left completely untouched, then stamp the same Courier New / 10pt formatting\n//    onto the new run.\nconst tracking = context.document.body.search(\n  \"Use project tracking tool \" + EN_DASH + \" Redmine, Jira\",\n  { matchCase: true }\n);\ntracking.load(\"text\");\nawait context.sync();\nif (tracking.items.length > 0) {\n  const trackingPara = tracking.items[0].paragraphs.getFirst();\n  const newRun = trackingPara.insertText(\".\", Word.InsertLocation.end);\n  newRun.font.name = \"Courier New\";\n  newRun.font.size = 10;\n  await context.sync();\n}\n\n// 4) Append \".\" right after the FIRST \"Experience in QA/QC workflow\" (same reasoning\n//    as step 3 above).\nconst qaqc = context.document.body.search(\"Experience in QA/QC workflow\", { matchCase: true });\nqaqc.load(\"text\");\nawait context.sync();\nif (qaqc.items.length > 0) {\n  const qaqcPara = qaqc.items[0].paragraphs.getFirst();\n  const newRun2 = qaqcPara.insertText(\".\", Word.InsertLocation.end);\n  newRun2.font.name = \"Courier New\";\n  newRun2.font.size = 10;\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Resume formatting fixes:\n# 1. \"2 years of experience in C or C++\" -> \"2 years of experience in C/C++\"\n# 2. \"Full Phase \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\" -> \"Full Phases \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\"\n# 3. First \"Use project tracking tool \u2013 Redmine, Jira\" gets a trailing \".\"\n# 4. First \"Experience in QA/QC workflow\" gets a trailing \".\"\n\n$d = $word.ActiveDocument\n\n# Word Find/Replace constants (not predefined by the host, so spelled out):\n#   wdReplaceNone = 0, wdReplaceOne = 1, wdReplaceAll = 2\n#   wdFindContinue = 1\n\n# 1) \"C or C++\" -> \"C/C++\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\"2 years of experience in C or C++\", $false, $false, $false, $false, $false, $true, 1, $false, \"2 years of experience in C/C++\", 2)\n\n# 2) \"Full Phase\" -> \"Full Phases\" (only instance of the project-name line)\n$rng2 = $d.Content\n$rng2.Find.Execute(\"Full Phase \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\", $false, $false, $false, $false, $false, $true, 1, $false, \"Full Phases \u2013 MCAL drivers \u2013 AUTOSAR (RVC)\", 2)\n\n# 3) Append \".\" after the first \"Use project tracking tool \u2013 Redmine, Jira\"\n$rng3 = $d.Content\n$rng3.Find.Execute(\"Use project tracking tool \u2013 Redmine, Jira\", $false, $false, $false, $false, $false, $true, 1, $false, \"Use project tracking tool \u2013 Redmine, Jira.\", 1)\n\n# 4) Append \".\" after the first \"Experience in QA/QC workflow\"\n$rng4 = $d.Content\n$rng4.Find.Execute(\"Experience in QA/QC workflow\", $false, $false, $false, $false, $false, $true, 1, $false, \"Experience in QA/QC workflow.\", 1)\n"}
